$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "MODEL_CONDITION" header text to "MODELCONDITION" (shared string used by E1).
$ws.Cells.Item(1, 5).Value = "MODELCONDITION"

# Remove column A entirely; the remaining columns (old B:F) shift left to become A:E.
$ws.Range("A1").EntireColumn.Delete()
